# Aggiornamento stato task Burndown chart
# Marks every "Non iniziata" (Not started) task status as "Completata" (Completed)
# on the "Task Sprint 1" sheet, and leaves that sheet active/selected (matching
# where the author was working when they saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task Sprint 1")

# Every "Stato" (status) cell in column E currently reads "Non iniziata";
# flip them all to "Completata".
$statusCells = @("E3","E5","E6","E7","E8","E10","E12","E13","E14","E16","E17","E18","E20","E21","E22","E24","E25")
foreach ($cell in $statusCells) {
    $ws.Range($cell).Value = "Completata"
}

# B24 previously carried a redundant one-off style; align it with the other
# empty "Story points" cells in its column (vertical-top alignment, style used
# throughout column B).
$ws.Range("B24").VerticalAlignment = -4160

# Leave "Task Sprint 1" as the active sheet/tab, with its last-worked cell
# selected.
$ws.Activate()
$ws.Range("J20").Select()
